$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing "Sprungschanze" row (row 6) for the two
# new "Stadion" values, pushing the remaining rows down.
$ws.Rows.Item(6).Insert()
$ws.Rows.Item(6).Insert()

$ws.Range("A6").Value = "Stadion, überdacht"
$ws.Range("B6").Value = 1441
$ws.Range("C6").Value = "http://inspire.ec.europa.eu/codelist/CurrentUseValue/publicServices"

$ws.Range("A7").Value = "Stadion, nicht überdacht"
$ws.Range("B7").Value = 1442
$ws.Range("C7").Value = "http://inspire.ec.europa.eu/codelist/CurrentUseValue/publicServices"

# Append a new row at the bottom of the table for "Wassersportanlage".
$ws.Range("A10").Value = "Wassersportanlage"
$ws.Range("B10").Value = 1650
$ws.Range("C10").Value = "http://inspire.ec.europa.eu/codelist/CurrentUseValue/publicServices"

# Match the style used by the other "value" rows in columns A-C.
$ws.Range("C10").Style = $ws.Range("C9").Style

# Remove the old multi-cell hyperlink (C3:C7) that duplicated the display text;
# only the single hyperlink on C2 should remain.
for ($i = $ws.Hyperlinks.Count; $i -ge 1; $i--) {
    $hl = $ws.Hyperlinks.Item($i)
    if ($hl.Range.Address() -ne '$C$2') {
        $hl.Delete()
    }
}

$ws.Range("C10").Select()
